$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data per Dec 1 2024 refresh

$ws.Cells.Item(2, 4).Value = "'97.518.77"
$ws.Cells.Item(2, 5).Value = "'  +1.06%  "

$ws.Cells.Item(3, 4).Value = "'3.717.38"
$ws.Cells.Item(3, 5).Value = "'  +0.04%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "'  -0.05%  "

$ws.Cells.Item(5, 4).Value = "'2.26"
$ws.Cells.Item(5, 5).Value = "'  +16.71%  "

$ws.Cells.Item(6, 4).Value = "'238.41"
$ws.Cells.Item(6, 5).Value = "'  +0.00%  "

$ws.Cells.Item(7, 5).Value = "'  +0.29%  "

$ws.Cells.Item(8, 4).Value = "'0.444"
$ws.Cells.Item(8, 5).Value = "'  +5.03%  "

$ws.Cells.Item(9, 4).Value = "'1.16"
$ws.Cells.Item(9, 5).Value = "'  +7.14%  "

$ws.Cells.Item(10, 5).Value = "'  -0.04%  "

$ws.Cells.Item(11, 4).Value = "'3.713.50"
$ws.Cells.Item(11, 5).Value = "'  +0.03%  "

$ws.Cells.Item(12, 4).Value = "'45.61"
$ws.Cells.Item(12, 5).Value = "'  +1.17%  "

$ws.Cells.Item(13, 4).Value = "'0.0000309"
$ws.Cells.Item(13, 5).Value = "'  +15.38%  "

$ws.Cells.Item(15, 4).Value = "'6.84"
$ws.Cells.Item(15, 5).Value = "'  +0.11%  "

$ws.Cells.Item(16, 4).Value = "'4.412.70"
$ws.Cells.Item(16, 5).Value = "'  +0.11%  "

$ws.Cells.Item(17, 4).Value = "'97.211.97"
$ws.Cells.Item(17, 5).Value = "'  +0.95%  "

$ws.Cells.Item(18, 4).Value = "'8.90"
$ws.Cells.Item(18, 5).Value = "'  -1.09%  "

$ws.Cells.Item(19, 4).Value = "'3.709.13"
$ws.Cells.Item(19, 5).Value = "'  -0.53%  "

$ws.Cells.Item(20, 4).Value = "'13.08"
$ws.Cells.Item(20, 5).Value = "'  +2.08%  "

$ws.Cells.Item(21, 4).Value = "'18.99"
$ws.Cells.Item(21, 5).Value = "'  -0.57%  "

$ws.Cells.Item(22, 4).Value = "'0.551"
$ws.Cells.Item(22, 5).Value = "'  +4.09%  "

$ws.Cells.Item(23, 4).Value = "'531.33"
$ws.Cells.Item(23, 5).Value = "'  +1.04%  "

$ws.Cells.Item(24, 5).Value = "'  -0.31%  "

$ws.Cells.Item(25, 5).Value = "'  +10.41%  "

$ws.Cells.Item(26, 4).Value = "'119.86"
$ws.Cells.Item(26, 5).Value = "'  +16.75%  "

$ws.Cells.Item(27, 4).Value = "'6.92"
$ws.Cells.Item(27, 5).Value = "'  -1.46%  "

$ws.Cells.Item(28, 4).Value = "'0.212"
$ws.Cells.Item(28, 5).Value = "'  +25.74%  "

$ws.Cells.Item(29, 4).Value = "'13.47"
$ws.Cells.Item(29, 5).Value = "'  +0.11%  "

$ws.Cells.Item(30, 4).Value = "'12.86"
$ws.Cells.Item(30, 5).Value = "'  +2.94%  "

$ws.Cells.Item(31, 4).Value = "'3.05"
$ws.Cells.Item(31, 5).Value = "'  -0.63%  "

$ws.Cells.Item(32, 5).Value = "'  -0.06%  "

$ws.Cells.Item(33, 4).Value = "'0.192"
$ws.Cells.Item(33, 5).Value = "'  +2.69%  "

$ws.Cells.Item(34, 2).Value = "'Binance-PegBSC-USD"
$ws.Cells.Item(34, 3).Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(34, 4).Value = "'1.02"
$ws.Cells.Item(34, 5).Value = "'  +2.23%  "

$ws.Cells.Item(35, 2).Value = "'Fetch.AI"
$ws.Cells.Item(35, 3).Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(35, 4).Value = "'1.83"
$ws.Cells.Item(35, 5).Value = "'  -3.64%  "

$ws.Cells.Item(36, 2).Value = "'EthereumClassic"
$ws.Cells.Item(36, 3).Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(36, 4).Value = "'33.16"
$ws.Cells.Item(36, 5).Value = "'  +1.03%  "

$ws.Cells.Item(37, 4).Value = "'0.605"
$ws.Cells.Item(37, 5).Value = "'  +0.85%  "

$ws.Cells.Item(38, 4).Value = "'640.55"

$ws.Cells.Item(39, 4).Value = "'8.73"
$ws.Cells.Item(39, 5).Value = "'  -1.75%  "

$ws.Cells.Item(40, 5).Value = "'  +0.01%  "

$ws.Cells.Item(41, 5).Value = "'  +4.58%  "

$ws.Cells.Item(42, 4).Value = "'6.92"
$ws.Cells.Item(42, 5).Value = "'  -3.17%  "

$ws.Cells.Item(43, 2).Value = "'Algorand"
$ws.Cells.Item(43, 3).Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(43, 4).Value = "'0.494"
$ws.Cells.Item(43, 5).Value = "'  +13.21%  "

$ws.Cells.Item(44, 2).Value = "'EnergySwap"
$ws.Cells.Item(44, 3).Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).Value = "'40.82"
$ws.Cells.Item(44, 5).Value = "'  +1.61%  "

$ws.Cells.Item(45, 5).Value = "'  +2.43%  "

$ws.Cells.Item(46, 4).Value = "'0.969"
$ws.Cells.Item(46, 5).Value = "'  -1.09%  "

$ws.Cells.Item(47, 5).Value = "'  +0.35%  "

$ws.Cells.Item(48, 4).Value = "'2.40"
$ws.Cells.Item(48, 5).Value = "'  +3.42%  "

$ws.Cells.Item(49, 4).Value = "'9.01"
$ws.Cells.Item(49, 5).Value = "'  +4.44%  "

$ws.Cells.Item(50, 4).Value = "'23.64"
$ws.Cells.Item(50, 5).Value = "'  +0.14%  "

$ws.Cells.Item(51, 4).Value = "'3.42"
$ws.Cells.Item(51, 5).Value = "'  +6.77%  "
